$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.597.83'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '3.509.68'
$ws.Range("E3").Value = '  -0.16%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '609.20'
$ws.Range("E5").Value = '  +0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.70'
$ws.Range("E6").Value = '  +0.96%  '
$ws.Range("D7").Value = '3.510.33'
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.487'
$ws.Range("E9").Value = '  +1.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.143'
$ws.Range("E10").Value = '  +2.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.62'
$ws.Range("E11").Value = '  +7.93%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.433'
$ws.Range("E12").Value = '  +2.15%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '32.53'
$ws.Range("E13").Value = '  +2.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000216'
$ws.Range("E14").Value = '  -1.37%  '
$ws.Range("D15").Value = '4.100.59'
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").Value = '3.508.93'
$ws.Range("E16").Value = '  -0.06%  '
$ws.Range("D17").Value = '67.804.43'
$ws.Range("E17").Value = '  +0.58%  '
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.56'
$ws.Range("E19").Value = '  +2.47%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.59'
$ws.Range("E20").Value = '  +2.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.86'
$ws.Range("E21").Value = '  +7.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '448.90'
$ws.Range("E22").Value = '  +0.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.633'
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.08'
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000128'
$ws.Range("E25").Value = '  -0.70%  '
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '3.643.69'
$ws.Range("E26").Value = '  -0.20%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.86'
$ws.Range("E28").Value = '  +6.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.12'
$ws.Range("E29").Value = '  -1.08%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.52'
$ws.Range("E30").Value = '  +0.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.65'
$ws.Range("E31").Value = '  +5.49%  '
$ws.Range("E32").Value = '  +3.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.76'
$ws.Range("E34").Value = '  +0.26%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.18'
$ws.Range("E35").Value = '  +1.06%  '
$ws.Range("E36").Value = '  +2.02%  '
$ws.Range("D37").Value = '3.500.99'
$ws.Range("E37").Value = '  -0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.03'
$ws.Range("E38").Value = '  -0.03%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.32'
$ws.Range("E40").Value = '  +6.49%  '
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("B42").Value = 'Hedera'
$ws.Range("C42").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0898'
$ws.Range("E42").Value = '  +2.77%  '
$ws.Range("B43").Value = 'Monero'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '173.82'
$ws.Range("E43").Value = '  -1.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.48'
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '30.34'
$ws.Range("E45").Value = '  +12.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.880'
$ws.Range("E46").Value = '  +0.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.59'
$ws.Range("E47").Value = '  +2.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.31'
$ws.Range("E48").Value = '  +3.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.67'
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("E50").Value = '  -2.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.254'
$ws.Range("E51").Value = '  +3.14%  '
